$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "63.840.10"
$ws.Range("E2").Value2 = "  +0.00%  "
$ws.Range("D3").Value2 = "2.624.04"
$ws.Range("E3").Value2 = "  -0.04%  "
$ws.Range("E4").Value2 = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "595.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value2 = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "151.15"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value2 = "  +0.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.586"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value2 = "  -0.19%  "
$ws.Range("E9").Value2 = "  +4.19%  "
$ws.Range("B10").Value2 = "Toncoin"
$ws.Range("C10").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "5.81"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value2 = "  +3.55%  "
$ws.Range("B11").Value2 = "Cardano"
$ws.Range("C11").Value2 = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.394"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value2 = "  +3.29%  "
$ws.Range("E12").Value2 = "  +1.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "27.97"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value2 = "  +1.48%  "
$ws.Range("D14").Value2 = "3.092.25"
$ws.Range("E14").Value2 = "  -0.09%  "
$ws.Range("D15").Value2 = "63.685.46"
$ws.Range("E15").Value2 = "  +0.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "0.0000165"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value2 = "  +10.79%  "
$ws.Range("D17").Value2 = "2.657.25"
$ws.Range("E17").Value2 = "  +1.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "12.24"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value2 = "  +0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "4.80"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value2 = "  +4.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "348.32"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value2 = "  -0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "7.02"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value2 = "  +1.88%  "
$ws.Range("E22").Value2 = "  +0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "67.47"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value2 = "  +2.00%  "
$ws.Range("E24").Value2 = "  -2.60%  "
$ws.Range("B25").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C25").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "9.22"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value2 = "  +0.45%  "
$ws.Range("B26").Value2 = "Fetch.AI"
$ws.Range("C26").Value2 = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "1.68"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value2 = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "8.42"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value2 = "  +3.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "550.45"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value2 = "  +1.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "0.163"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value2 = "  -0.90%  "
$ws.Range("E30").Value2 = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "2.07"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value2 = "  +1.55%  "
$ws.Range("D32").Value2 = "0.0₃0892"
$ws.Range("E32").Value2 = "  +5.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "1.80"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value2 = "  +3.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "5.42"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value2 = "  +4.00%  "
$ws.Range("E35").Value2 = "  +2.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "164.51"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value2 = "  -2.45%  "
$ws.Range("E37").Value2 = "  +2.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "19.81"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value2 = "  +2.26%  "
$ws.Range("E40").Value2 = "  -0.07%  "
$ws.Range("E41").Value2 = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "167.68"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value2 = "  -1.50%  "
$ws.Range("E43").Value2 = "  +4.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "23.62"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value2 = "  +10.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "0.0586"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value2 = "  -1.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "2.18"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value2 = "  +10.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.637"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value2 = "  +1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "0.0254"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value2 = "  +3.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0969"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value2 = "  +0.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "19.26"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value2 = "  +0.44%  "
$ws.Range("E51").Value2 = "  +18.05%  "
